$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the Price column keeps being stored as text (it contains values
# like "56.845.10" that are not valid numbers, mixed with values like
# "490.37" that Excel would otherwise auto-convert to a number).
$ws.Range("D2:D51").NumberFormat = "@"

# Simple price / volume updates (no coin name / link changes)
$ws.Range("D2").Value = "56.845.10"
$ws.Range("E2").Value = "  +5.81%  "

$ws.Range("D3").Value = "2.484.48"
$ws.Range("E3").Value = "  +3.54%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "490.37"
$ws.Range("E5").Value = "  +6.23%  "

$ws.Range("D6").Value = "146.55"
$ws.Range("E6").Value = "  +11.95%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  +6.42%  "

$ws.Range("D9").Value = "2.501.55"
$ws.Range("E9").Value = "  +3.67%  "

$ws.Range("D10").Value = "5.79"
$ws.Range("E10").Value = "  +9.77%  "

$ws.Range("D11").Value = "0.0976"
$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  +6.04%  "

$ws.Range("D13").Value = "0.124"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "2.922.14"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").Value = "56.734.16"
$ws.Range("E15").Value = "  +5.49%  "

$ws.Range("D16").Value = "21.26"
$ws.Range("E16").Value = "  +8.12%  "

$ws.Range("E17").Value = "  +5.15%  "

$ws.Range("D18").Value = "2.504.62"
$ws.Range("E18").Value = "  +3.52%  "

$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  +9.95%  "

$ws.Range("D20").Value = "10.24"
$ws.Range("E20").Value = "  +9.99%  "

$ws.Range("D21").Value = "320.07"
$ws.Range("E21").Value = "  +4.64%  "

$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("E23").Value = "  +9.30%  "

$ws.Range("D24").Value = "59.06"
$ws.Range("E24").Value = "  +5.37%  "

$ws.Range("D25").Value = "0.413"
$ws.Range("E25").Value = "  +8.03%  "

$ws.Range("E26").Value = "  +8.44%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "2.594.69"
$ws.Range("E28").Value = "  +2.06%  "

$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  +8.17%  "

$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  +11.11%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").Value = "149.35"
$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").Value = "18.24"
$ws.Range("E33").Value = "  +3.53%  "

$ws.Range("E34").Value = "  +6.14%  "

$ws.Range("D35").Value = "5.21"
$ws.Range("E35").Value = "  +4.78%  "

$ws.Range("E36").Value = "  +9.14%  "

$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  +6.59%  "

$ws.Range("D38").Value = "0.864"
$ws.Range("E38").Value = "  +9.18%  "

$ws.Range("D39").Value = "34.24"
$ws.Range("E39").Value = "  +3.81%  "

$ws.Range("E40").Value = "  +8.27%  "

# Rows 41 and 42 swap places (Hedera <-> Mantle) with updated price/volume
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "0.614"
$ws.Range("E41").Value = "  +3.79%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0559"
$ws.Range("E42").Value = "  +7.11%  "

$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  +8.72%  "

$ws.Range("D45").Value = "4.82"
$ws.Range("E45").Value = "  +14.69%  "

$ws.Range("D46").Value = "259.61"
$ws.Range("E46").Value = "  +18.69%  "

$ws.Range("D47").Value = "0.0230"
$ws.Range("E47").Value = "  +6.07%  "

$ws.Range("D48").Value = "0.0918"
$ws.Range("E48").Value = "  +6.14%  "

$ws.Range("D49").Value = "10.21"
$ws.Range("E49").Value = "  +0.28%  "

# Rows 50 and 51 swap places (Maker <-> EnergySwap) with updated price/volume
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "17.71"
$ws.Range("E50").Value = "  +7.83%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.898.53"
$ws.Range("E51").Value = "  -1.63%  "
